# Daily update at 8 AM UTC
# Appends the day's new snapshot row to the bottom of the "Wins Over Time"
# table and moves the "most recent entry" date-only formatting down to the
# newly appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row in column A (the "Day" column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# The row that was previously "last" used a date-only number format to mark
# it as the latest entry; now that a newer row follows it, it reverts to the
# standard date/time format shared by the rest of the table.
$ws.Range("A" + $lastRow).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$newRow = $lastRow + 1

# Use Value2 to get the raw date serial number (Value would return a
# formatted string since the cell has a date number format applied).
$prevDate = $ws.Range("A" + $lastRow).Value2

# New daily snapshot: one day after the previous entry, with Chase, Bryce,
# and Zach all tied at the same win count.
$ws.Cells.Item($newRow, 1).Value = $prevDate + 1
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 2).Value = 177
$ws.Cells.Item($newRow, 3).Value = 177
$ws.Cells.Item($newRow, 4).Value = 177

Write-Host "Appended daily snapshot to row $newRow"
